# Automatische test-sync: 2025-08-04 20:53:50
#
# 1) Logs sheet: append row 25 (Testmail #13 - VentiQ-250 datasheet request)
#    and extend the conditional-formatting ranges that referenced the old
#    last row (24) so they also cover the new last row (25).
# 2) Dashboard sheet: swap the category labels that used to sit in rows 2/3
#    and append a new row 7 for the "Documentatie / Datasheets" category.
# 3) Chart1 on the Dashboard sheet: extend its category/value series
#    references from row 6 to row 7 so the new category is plotted.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A25").Value = "Kun je mij de datasheet van de VentiQ-250 sturen?"
$logs.Range("B25").Value = "mailmind.test@zohomail.eu"
$logs.Range("C25").Value = "Testmail #13: Kun je mij de datasheet van de VentiQ-250 sturen?"
$logs.Range("D25").Value = "Documentatie / Datasheets"
$logs.Range("E25").Value = "Bedankt, we hebben dit doorgestuurd naar documentatie@bedrijf.nl."
$logs.Range("F25").Value = "2025-08-04 20:53:41"
$logs.Range("G25").Value = "Ja"
$logs.Range("H25").Value = "Ja"
$logs.Range("I25").Value = "Nee"
$logs.Range("J25").Value = "Nee"

# Extend every conditional-formatting block that used to stop at row 24 so
# it covers the freshly added row 25 as well (sqref D2:D24 -> D2:D25, etc.)
function Extend-ConditionalFormatting($oldRange, $newRange) {
    $fcs = $logs.Range($oldRange).FormatConditions
    $count = $fcs.Count()
    for ($i = 1; $i -le $count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newRange))
    }
}

Extend-ConditionalFormatting "D2:D24" "D2:D25"
Extend-ConditionalFormatting "G2:G24" "G2:G25"
Extend-ConditionalFormatting "H2:H24" "H2:H25"
Extend-ConditionalFormatting "I2:I24" "I2:I25"
Extend-ConditionalFormatting "J2:J24" "J2:J25"

# ---------------------------------------------------------------------
# 2) Dashboard sheet
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

# Rows 2 and 3 swap their category labels (counts stay 6 / 6, unchanged).
$dash.Range("A2").Value = "Inkoop / Bestellingen"
$dash.Range("A3").Value = "Planning / Afspraak"

# New row 7 for the "Documentatie / Datasheets" category.
$dash.Range("A7").Value = "Documentatie / Datasheets"
$dash.Range("B7").Value = 1

# ---------------------------------------------------------------------
# 3) Chart1: extend category/value series ranges from row 6 to row 7
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$7,Dashboard!`$B`$2:`$B`$7,1)"
